$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.24%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '13'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.95%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '13'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.129'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.70%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '13'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07832'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.37%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '13'

# Row 6
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.400'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.38%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '13'

# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.283'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.92%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '13'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.878'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.29%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '13'

# Row 9
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.999'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.92%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '13'

# Row 10
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9225'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.30%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '13'

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1075'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-9.00%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '13'

# Row 12
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1896'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.15%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '13'

# Row 13
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08828'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-6.67%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '13'

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03321'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-3.83%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '13'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09594'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.32%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '13'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001390'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.97%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '13'

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005823'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.09%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '13'

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.397'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.85%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '13'

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.41%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '13'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.338'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '20.55%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '13'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1299'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.47%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '13'

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-6.97%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '13'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04351'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.44%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '13'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001201'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.38%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '13'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004263'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.18%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '13'

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '7.75%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '13'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002907'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '13'

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '13'

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '13'

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '13'

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '13'

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '13'

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '13'

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '13'

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '13'

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '13'

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '13'

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '13'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02151'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.96%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '13'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05016'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.65%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '13'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007579'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.17%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '13'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1352'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.47%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '13'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008522'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.64%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '13'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002012'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.89%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '13'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008123'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-5.69%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '13'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006549'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.30%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '13'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000752'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.24%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '13'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003299'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '12.43%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '13'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001446'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '20.59%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '13'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002105'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.24%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '13'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002005'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.24%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '13'
